$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.979584693908691
$ws.Range("B1").Value = 4.402648448944092
$ws.Range("C1").Value = 2.11156964302063
$ws.Range("D1").Value = 1.604492545127869
$ws.Range("E1").Value = 1.427490472793579
